$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append the new row of data (row 17) that was recorded on 2025-04-29.
# The date column holds free-form text like "MM/DD/YYYY" (not a real
# Excel date), so force the cell to text first to avoid Excel
# auto-converting the string into a date serial number, then restore
# the default (Normal) style so no extra number formatting is applied.
$ws.Range("A17").NumberFormat = "@"
$ws.Range("A17").Value = "04/29/2025"
$ws.Range("A17").Style = "Normal"

$ws.Range("B17").Value = 535.3930000000037
$ws.Range("C17").Value = 0.09338934203472898
$ws.Range("D17").Value = 50
